# Auto-generated: apply cell value updates from the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.177.29"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.423.52"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'553.88"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'137.08"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").Value = "'5.72"
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").Value = "'24.86"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Value = "2.857.66"
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "60.103.83"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Value = "'0.0000138"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "2.415.85"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'11.26"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "'4.50"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").Value = "'327.19"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D23").Value = "'65.31"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +4.18%  "
$ws.Range("D25").Value = "'8.65"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  +4.85%  "
$ws.Range("D28").Value = "0.0₃0774"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'170.51"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").Value = "'6.11"
$ws.Range("E31").Value = "  -2.45%  "
$ws.Range("E32").Value = "  -3.20%  "
$ws.Range("D33").Value = "'1.06"
$ws.Range("D34").Value = "'18.54"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("E35").Value = "  +2.58%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.22"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'325.81"
$ws.Range("E39").Value = "  +3.73%  "
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "'145.56"
$ws.Range("E41").Value = "  +4.52%  "
$ws.Range("D42").Value = "'3.65"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").Value = "'0.0963"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "'19.84"
$ws.Range("E44").Value = "  +1.64%  "
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "'0.0223"
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").Value = "'4.66"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("E51").Value = "  -0.70%  "
